# Weekly refresh of the Berenjena (Hortaliza) price series:
# a new daily observation is inserted as the first data row (row 10),
# pushing the existing rows 10-28 down to 11-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 10:28 down to 11:29 and open up a blank row 10.
$ws.Range("A10").EntireRow.Insert()

# Populate the newly inserted row 10 with the latest observation.
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44690
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 100112001
$ws.Range("G10").Value = "Berenjena"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 60
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = 10000
$ws.Range("N10").Value = "$/caja 60 unidades"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 167
$ws.Range("Q10").Value = 60
$ws.Range("R10").Value = "Hortaliza"
